$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing data rows down
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new credential record
$ws.Range("A2").Value = "Walter"
$ws.Range("B2").Value = "Newport"
$ws.Range("C2").Value = "walteru6_newportgw@outlook.com"
$ws.Range("D2").Value = "wba18Hq1BbNr5xd"
$ws.Range("E2").Value = "81.28.96.40:4004"
$ws.Range("F2").Value = "PJ5C8sm37i4b"
$ws.Range("G2").Value = "56gRMx51KSrg"

# Update the active selection to match the saved view state
$null = $ws.Range("D8").Select()
